$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13036
$ws.Range("C2").Value = 1262
$ws.Range("B3").Value = 704
$ws.Range("C3").Value = 4197
